$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.405
$ws.Range("D7").Value = -7.768000000000001
$ws.Range("A9").Value = -21.658
$ws.Range("D12").Value = -7.456
$ws.Range("A13").Value = -22.219
$ws.Range("D14").Value = -7.813
$ws.Range("E15").Value = 15.94
$ws.Range("A16").Value = -22.013
$ws.Range("A18").Value = -21.93
$ws.Range("D19").Value = -7.851999999999999
$ws.Range("A20").Value = -20.1
$ws.Range("A26").Value = -21.59
$ws.Range("D26").Value = -7.547
$ws.Range("A27").Value = -22.004
$ws.Range("D27").Value = -8.228999999999999
$ws.Range("E28").Value = 17.041
$ws.Range("A29").Value = -21.135
$ws.Range("D29").Value = -7.481999999999999
$ws.Range("E33").Value = 17.093
$ws.Range("A35").Value = -19.982
$ws.Range("E35").Value = 16.604
$ws.Range("A36").Value = -20.339
$ws.Range("D37").Value = -7.822
$ws.Range("D38").Value = -7.795999999999999
$ws.Range("E38").Value = 16.796
$ws.Range("E43").Value = 17.052
$ws.Range("E44").Value = 16.723
$ws.Range("A45").Value = -21.453
$ws.Range("E45").Value = 16.526
$ws.Range("D47").Value = -7.325
$ws.Range("E47").Value = 16.854
$ws.Range("D51").Value = -8.407
$ws.Range("E51").Value = 16.726
$ws.Range("D52").Value = -8.1
$ws.Range("E54").Value = 16.655
$ws.Range("A55").Value = -22.149
$ws.Range("D55").Value = -8.041999999999998
$ws.Range("A57").Value = -22.125
$ws.Range("E57").Value = 16.35
$ws.Range("E62").Value = 16.207
$ws.Range("E63").Value = 17.689
$ws.Range("E67").Value = 17.018
$ws.Range("A69").Value = -21.498
$ws.Range("D69").Value = -7.56
$ws.Range("D70").Value = -7.531000000000001
$ws.Range("E70").Value = 17.585
$ws.Range("A76").Value = -20.047
$ws.Range("D76").Value = -7.743
$ws.Range("A78").Value = -20.086
$ws.Range("D81").Value = -7.962000000000001
$ws.Range("E81").Value = 17.223
$ws.Range("A82").Value = -22.018
$ws.Range("A83").Value = -21.872
$ws.Range("D83").Value = -8.270999999999999
$ws.Range("E88").Value = 16.447
$ws.Range("A93").Value = -21.449
$ws.Range("D94").Value = -7.456
$ws.Range("E96").Value = 16.725
$ws.Range("A97").Value = -22.136
$ws.Range("E99").Value = 16.775
$ws.Range("D100").Value = -8.215999999999999
$ws.Range("D102").Value = -7.764
